$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two trailing rows (record_type/recordtype/other/10 and
# Event Type/eventtype/dexcom/10 at rows 13 & 14). Row 14 goes first so the
# row index of the remaining row to delete does not shift.
$ws.Rows.Item(14).Delete()
$ws.Rows.Item(13).Delete()

# Row 5: Record Type now maps to the "eventtype" variable instead of scan_yn
$ws.Range("B5").Value = "eventtype"

# Row 6: previously "Scan Glucose(mmol/L)" / scanglucose / libre / 14,
# now becomes "Source Device ID" / deviceid / dexcom / 10 (moved up from
# the old row 9 position).
$ws.Range("A6").Value = "Source Device ID"
$ws.Range("B6").Value = "deviceid"
$ws.Range("C6").Value = "dexcom"
$ws.Range("D6").Value = 10

# Give A6 the same cell style already used by A7/A8/A9/A12 (bold-ish
# Helvetica Neue font) by copying the format from A7.
$ws.Range("A7").Copy()
$ws.Range("A6").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row 9: now holds "Event Type" / eventtype / dexcom / 10 (previously held
# "Source Device ID" / deviceid / dexcom / 10, which moved to row 6).
$ws.Range("A9").Value = "Event Type"
$ws.Range("B9").Value = "eventtype"
$ws.Range("C9").Value = "dexcom"
$ws.Range("D9").Value = 10

# Update the saved selection to match the target sheet view.
$ws.Range("F17").Select()
